$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Flip the existing funct4=0/1 truth-table rows (12,14,16,18): swap O/P ---
$ws.Range("O12").Value = "'1"
$ws.Range("P12").Value = "'0"

$ws.Range("O14").Value = "'1"
$ws.Range("P14").Value = "'0"

$ws.Range("O16").Value = "'1"
$ws.Range("P16").Value = "'0"

$ws.Range("O18").Value = "'1"
$ws.Range("P18").Value = "'0"

# --- Append the new OPSET funct4 = 0101 entry (rows 19-20) ---
$ws.Range("N19").Value = "'0101"
$ws.Range("O19").Value = "'0"
$ws.Range("P19").Value = "'0"
$ws.Range("Q19").Value = "I"

$ws.Range("O20").Value = "'1"
$ws.Range("P20").Value = "'0"
$ws.Range("Q20").Value = "X"
